$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.931376
$ws.Range("H2").Value = 8.794128
$ws.Range("I2").Value = 0.1253619302628033
$ws.Range("J2").Value = 0.1253619302628033
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 126.660156
$ws.Range("N2").Value = 379.980468
$ws.Range("O2").Value = 0.9724734733029194
$ws.Range("P2").Value = 0.9724734733029196
$ws.Range("Q2").Value = 371.288541454656
$ws.Range("R2").Value = 3341.596873091904
$ws.Range("S2").Value = 0.1219111517426267
$ws.Range("T2").Value = 0.1219111517426267

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.931376
$ws.Range("H3").Value = 8.794128
$ws.Range("I3").Value = 0.1253619302628033
$ws.Range("J3").Value = 0.1253619302628033
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.249583
$ws.Range("N3").Value = 0.7487489999999999
$ws.Range("O3").Value = 0.001916252549754972
$ws.Range("P3").Value = 0.001916252549754972
$ws.Range("Q3").Value = 0.7316216162079999
$ws.Range("R3").Value = 6.584594545871999
$ws.Range("S3").Value = 0.0002402251185083019
$ws.Range("T3").Value = 0.0002402251185083019

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.931376
$ws.Range("H4").Value = 8.794128
$ws.Range("I4").Value = 0.1253619302628033
$ws.Range("J4").Value = 0.1253619302628033
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.048924
$ws.Range("N4").Value = 3.146772
$ws.Range("O4").Value = 0.008053446306435876
$ws.Range("P4").Value = 0.008053446306435876
$ws.Range("Q4").Value = 3.074790639424001
$ws.Range("R4").Value = 27.67311575481601
$ws.Range("S4").Value = 0.001009595574242645
$ws.Range("T4").Value = 0.001009595574242645

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.931376
$ws.Range("H5").Value = 8.794128
$ws.Range("I5").Value = 0.1253619302628033
$ws.Range("J5").Value = 0.1253619302628033
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.522703666666667
$ws.Range("N5").Value = 4.568111
$ws.Range("O5").Value = 0.01169103978945379
$ws.Range("P5").Value = 0.01169103978945379
$ws.Range("Q5").Value = 4.463616983578667
$ws.Range("R5").Value = 40.172552852208
$ws.Range("S5").Value = 0.001465611314785165
$ws.Range("T5").Value = 0.001465611314785165

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.931376
$ws.Range("H6").Value = 8.794128
$ws.Range("I6").Value = 0.1253619302628033
$ws.Range("J6").Value = 0.1253619302628033
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.7639916666666666
$ws.Range("N6").Value = 2.291975
$ws.Range("O6").Value = 0.005865788051435999
$ws.Range("P6").Value = 0.005865788051435999
$ws.Range("Q6").Value = 2.239546835866667
$ws.Range("R6").Value = 20.1559215228
$ws.Range("S6").Value = 0.0007353465126405048
$ws.Range("T6").Value = 0.0007353465126405047

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.220282666666667
$ws.Range("H7").Value = 15.660848
$ws.Range("I7").Value = 0.2232483010063491
$ws.Range("J7").Value = 0.2232483010063491
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 126.660156
$ws.Range("N7").Value = 379.980468
$ws.Range("O7").Value = 0.9724734733029194
$ws.Range("P7").Value = 0.9724734733029196
$ws.Range("Q7").Value = 661.201816924096
$ws.Range("R7").Value = 5950.816352316863
$ws.Range("S7").Value = 0.21710305068862
$ws.Range("T7").Value = 0.21710305068862

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.220282666666667
$ws.Range("H8").Value = 15.660848
$ws.Range("I8").Value = 0.2232483010063491
$ws.Range("J8").Value = 0.2232483010063491
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.249583
$ws.Range("N8").Value = 0.7487489999999999
$ws.Range("O8").Value = 0.001916252549754972
$ws.Range("P8").Value = 0.001916252549754972
$ws.Range("Q8").Value = 1.302893808794667
$ws.Range("R8").Value = 11.726044279152
$ws.Range("S8").Value = 0.000427800126031882
$ws.Range("T8").Value = 0.000427800126031882

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.220282666666667
$ws.Range("H9").Value = 15.660848
$ws.Range("I9").Value = 0.2232483010063491
$ws.Range("J9").Value = 0.2232483010063491
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.048924
$ws.Range("N9").Value = 3.146772
$ws.Range("O9").Value = 0.008053446306435876
$ws.Range("P9").Value = 0.008053446306435876
$ws.Range("Q9").Value = 5.475679775850668
$ws.Range("R9").Value = 49.28111798265601
$ws.Range("S9").Value = 0.001797918205157667
$ws.Range("T9").Value = 0.001797918205157667

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.220282666666667
$ws.Range("H10").Value = 15.660848
$ws.Range("I10").Value = 0.2232483010063491
$ws.Range("J10").Value = 0.2232483010063491
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.522703666666667
$ws.Range("N10").Value = 4.568111
$ws.Range("O10").Value = 0.01169103978945379
$ws.Range("P10").Value = 0.01169103978945379
$ws.Range("Q10").Value = 7.948943557569778
$ws.Range("R10").Value = 71.54049201812799
$ws.Range("S10").Value = 0.002610004769993184
$ws.Range("T10").Value = 0.002610004769993185

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 5.220282666666667
$ws.Range("H11").Value = 15.660848
$ws.Range("I11").Value = 0.2232483010063491
$ws.Range("J11").Value = 0.2232483010063491
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.7639916666666666
$ws.Range("N11").Value = 2.291975
$ws.Range("O11").Value = 0.005865788051435999
$ws.Range("P11").Value = 0.005865788051435999
$ws.Range("Q11").Value = 3.988252454977778
$ws.Range("R11").Value = 35.89427209479999
$ws.Range("S11").Value = 0.00130952721654643
$ws.Range("T11").Value = 0.00130952721654643

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6.882553333333334
$ws.Range("H12").Value = 20.64766
$ws.Range("I12").Value = 0.2943362335651782
$ws.Range("J12").Value = 0.2943362335651782
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 126.660156
$ws.Range("N12").Value = 379.980468
$ws.Range("O12").Value = 0.9724734733029194
$ws.Range("P12").Value = 0.9724734733029196
$ws.Range("Q12").Value = 871.7452788783199
$ws.Range("R12").Value = 7845.70750990488
$ws.Range("S12").Value = 0.2862341793740282
$ws.Range("T12").Value = 0.2862341793740282

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 6.882553333333334
$ws.Range("H13").Value = 20.64766
$ws.Range("I13").Value = 0.2943362335651782
$ws.Range("J13").Value = 0.2943362335651782
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.249583
$ws.Range("N13").Value = 0.7487489999999999
$ws.Range("O13").Value = 0.001916252549754972
$ws.Range("P13").Value = 0.001916252549754972
$ws.Range("Q13").Value = 1.717768308593333
$ws.Range("R13").Value = 15.45991477734
$ws.Range("S13").Value = 0.0005640225580545477
$ws.Range("T13").Value = 0.0005640225580545478

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.882553333333334
$ws.Range("H14").Value = 20.64766
$ws.Range("I14").Value = 0.2943362335651782
$ws.Range("J14").Value = 0.2943362335651782
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.048924
$ws.Range("N14").Value = 3.146772
$ws.Range("O14").Value = 0.008053446306435876
$ws.Range("P14").Value = 0.008053446306435876
$ws.Range("Q14").Value = 7.219275372613335
$ws.Range("R14").Value = 64.97347835352001
$ws.Range("S14").Value = 0.002370421053055732
$ws.Range("T14").Value = 0.002370421053055732

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.882553333333334
$ws.Range("H15").Value = 20.64766
$ws.Range("I15").Value = 0.2943362335651782
$ws.Range("J15").Value = 0.2943362335651782
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.522703666666667
$ws.Range("N15").Value = 4.568111
$ws.Range("O15").Value = 0.01169103978945379
$ws.Range("P15").Value = 0.01169103978945379
$ws.Range("Q15").Value = 10.48008919669556
$ws.Range("R15").Value = 94.32080277026002
$ws.Range("S15").Value = 0.003441096618088463
$ws.Range("T15").Value = 0.003441096618088464

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.882553333333334
$ws.Range("H16").Value = 20.64766
$ws.Range("I16").Value = 0.2943362335651782
$ws.Range("J16").Value = 0.2943362335651782
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.7639916666666666
$ws.Range("N16").Value = 2.291975
$ws.Range("O16").Value = 0.005865788051435999
$ws.Range("P16").Value = 0.005865788051435999
$ws.Range("Q16").Value = 5.258213392055556
$ws.Range("R16").Value = 47.3239205285
$ws.Range("S16").Value = 0.001726513961951298
$ws.Range("T16").Value = 0.001726513961951298

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 5.259310999999999
$ws.Range("H17").Value = 15.777933
$ws.Range("I17").Value = 0.2249173694580273
$ws.Range("J17").Value = 0.2249173694580273
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 126.660156
$ws.Range("N17").Value = 379.980468
$ws.Range("O17").Value = 0.9724734733029194
$ws.Range("P17").Value = 0.9724734733029196
$ws.Range("Q17").Value = 666.1451517125158
$ws.Range("R17").Value = 5995.306365412644
$ws.Range("S17").Value = 0.2187261754830038
$ws.Range("T17").Value = 0.2187261754830039

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 5.259310999999999
$ws.Range("H18").Value = 15.777933
$ws.Range("I18").Value = 0.2249173694580273
$ws.Range("J18").Value = 0.2249173694580273
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 0.249583
$ws.Range("N18").Value = 0.7487489999999999
$ws.Range("O18").Value = 0.001916252549754972
$ws.Range("P18").Value = 0.001916252549754972
$ws.Range("Q18").Value = 1.312634617313
$ws.Range("R18").Value = 11.813711555817
$ws.Range("S18").Value = 0.0004309984827081259
$ws.Range("T18").Value = 0.000430998482708126

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 5.259310999999999
$ws.Range("H19").Value = 15.777933
$ws.Range("I19").Value = 0.2249173694580273
$ws.Range("J19").Value = 0.2249173694580273
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 1.048924
$ws.Range("N19").Value = 3.146772
$ws.Range("O19").Value = 0.008053446306435876
$ws.Range("P19").Value = 0.008053446306435876
$ws.Range("Q19").Value = 5.516617531364
$ws.Range("R19").Value = 49.64955778227601
$ws.Range("S19").Value = 0.001811359958315023
$ws.Range("T19").Value = 0.001811359958315023

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 5.259310999999999
$ws.Range("H20").Value = 15.777933
$ws.Range("I20").Value = 0.2249173694580273
$ws.Range("J20").Value = 0.2249173694580273
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 1.522703666666667
$ws.Range("N20").Value = 4.568111
$ws.Range("O20").Value = 0.01169103978945379
$ws.Range("P20").Value = 0.01169103978945379
$ws.Range("Q20").Value = 8.008372143840333
$ws.Range("R20").Value = 72.075349294563
$ws.Range("S20").Value = 0.002629517915673076
$ws.Range("T20").Value = 0.002629517915673077

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 5.259310999999999
$ws.Range("H21").Value = 15.777933
$ws.Range("I21").Value = 0.2249173694580273
$ws.Range("J21").Value = 0.2249173694580273
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 0.7639916666666666
$ws.Range("N21").Value = 2.291975
$ws.Range("O21").Value = 0.005865788051435999
$ws.Range("P21").Value = 0.005865788051435999
$ws.Range("Q21").Value = 4.018069776408333
$ws.Range("R21").Value = 36.16262798767499
$ws.Range("S21").Value = 0.001319317618327313
$ws.Range("T21").Value = 0.001319317618327313

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 3.08978
$ws.Range("H22").Value = 9.26934
$ws.Range("I22").Value = 0.1321361657076419
$ws.Range("J22").Value = 0.1321361657076419
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 126.660156
$ws.Range("N22").Value = 379.980468
$ws.Range("O22").Value = 0.9724734733029194
$ws.Range("P22").Value = 0.9724734733029196
$ws.Range("Q22").Value = 391.3520168056799
$ws.Range("R22").Value = 3522.16815125112
$ws.Range("S22").Value = 0.1284989160146406
$ws.Range("T22").Value = 0.1284989160146406

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 3.08978
$ws.Range("H23").Value = 9.26934
$ws.Range("I23").Value = 0.1321361657076419
$ws.Range("J23").Value = 0.1321361657076419
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 0.249583
$ws.Range("N23").Value = 0.7487489999999999
$ws.Range("O23").Value = 0.001916252549754972
$ws.Range("P23").Value = 0.001916252549754972
$ws.Range("Q23").Value = 0.7711565617399998
$ws.Range("R23").Value = 6.940409055659999
$ws.Range("S23").Value = 0.0002532062644521142
$ws.Range("T23").Value = 0.0002532062644521143

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 3.08978
$ws.Range("H24").Value = 9.26934
$ws.Range("I24").Value = 0.1321361657076419
$ws.Range("J24").Value = 0.1321361657076419
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 1.048924
$ws.Range("N24").Value = 3.146772
$ws.Range("O24").Value = 0.008053446306435876
$ws.Range("P24").Value = 0.008053446306435876
$ws.Range("Q24").Value = 3.24094439672
$ws.Range("R24").Value = 29.16849957048
$ws.Range("S24").Value = 0.001064151515664807
$ws.Range("T24").Value = 0.001064151515664807

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 3.08978
$ws.Range("H25").Value = 9.26934
$ws.Range("I25").Value = 0.1321361657076419
$ws.Range("J25").Value = 0.1321361657076419
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 1.522703666666667
$ws.Range("N25").Value = 4.568111
$ws.Range("O25").Value = 0.01169103978945379
$ws.Range("P25").Value = 0.01169103978945379
$ws.Range("Q25").Value = 4.704819335193333
$ws.Range("R25").Value = 42.34337401674
$ws.Range("S25").Value = 0.001544809170913901
$ws.Range("T25").Value = 0.001544809170913901

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 3.08978
$ws.Range("H26").Value = 9.26934
$ws.Range("I26").Value = 0.1321361657076419
$ws.Range("J26").Value = 0.1321361657076419
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 0.7639916666666666
$ws.Range("N26").Value = 2.291975
$ws.Range("O26").Value = 0.005865788051435999
$ws.Range("P26").Value = 0.005865788051435999
$ws.Range("Q26").Value = 2.360566171833333
$ws.Range("R26").Value = 21.2450955465
$ws.Range("S26").Value = 0.000775082741970453
$ws.Range("T26").Value = 0.000775082741970453
